$d = $word.ActiveDocument
$d.Content.Find.Execute("425×7=", $true, $false, $false, $false, $false, $true, 1, $false, "283×8=", 2)
$d.Content.Find.Execute("213×7=", $true, $false, $false, $false, $false, $true, 1, $false, "496×9=", 2)
$d.Content.Find.Execute("555×5=", $true, $false, $false, $false, $false, $true, 1, $false, "568×9=", 2)
$d.Content.Find.Execute("956×2=", $true, $false, $false, $false, $false, $true, 1, $false, "851×6=", 2)
$d.Content.Find.Execute("940×5=", $true, $false, $false, $false, $false, $true, 1, $false, "996×5=", 2)
$d.Content.Find.Execute("766×4=", $true, $false, $false, $false, $false, $true, 1, $false, "843×9=", 2)
$d.Content.Find.Execute("108×7=", $true, $false, $false, $false, $false, $true, 1, $false, "371×9=", 2)
$d.Content.Find.Execute("286×8=", $true, $false, $false, $false, $false, $true, 1, $false, "895×8=", 2)
$d.Content.Find.Execute("148×9=", $true, $false, $false, $false, $false, $true, 1, $false, "693×2=", 2)
$d.Content.Find.Execute("651×2=", $true, $false, $false, $false, $false, $true, 1, $false, "297×7=", 2)
$d.Content.Find.Execute("504×7=", $true, $false, $false, $false, $false, $true, 1, $false, "942×3=", 2)
$d.Content.Find.Execute("472×2=", $true, $false, $false, $false, $false, $true, 1, $false, "796×6=", 2)
$d.Content.Find.Execute("988×3=", $true, $false, $false, $false, $false, $true, 1, $false, "854×6=", 2)
$d.Content.Find.Execute("131×5=", $true, $false, $false, $false, $false, $true, 1, $false, "446×3=", 2)
$d.Content.Find.Execute("521×6=", $true, $false, $false, $false, $false, $true, 1, $false, "736×8=", 2)
$d.Content.Find.Execute("798×6=", $true, $false, $false, $false, $false, $true, 1, $false, "227×7=", 2)
$d.Content.Find.Execute("988×6=", $true, $false, $false, $false, $false, $true, 1, $false, "282×5=", 2)
$d.Content.Find.Execute("625×3=", $true, $false, $false, $false, $false, $true, 1, $false, "708×2=", 2)
$d.Content.Find.Execute("186×6=", $true, $false, $false, $false, $false, $true, 1, $false, "654×3=", 2)
$d.Content.Find.Execute("404×4=", $true, $false, $false, $false, $false, $true, 1, $false, "923×8=", 2)
$d.Content.Find.Execute("780×7=", $true, $false, $false, $false, $false, $true, 1, $false, "922×3=", 2)
$d.Content.Find.Execute("903×9=", $true, $false, $false, $false, $false, $true, 1, $false, "357×2=", 2)
$d.Content.Find.Execute("838×9=", $true, $false, $false, $false, $false, $true, 1, $false, "222×6=", 2)
$d.Content.Find.Execute("259×3=", $true, $false, $false, $false, $false, $true, 1, $false, "481×4=", 2)
$d.Content.Find.Execute("287×4=", $true, $false, $false, $false, $false, $true, 1, $false, "521×7=", 2)
